# "Changed column from 'deprecated?' to 'state'"
#
# Column E was a boolean "Deprecated" flag (TRUE/FALSE, the first two rows
# even carried a literal =TRUE() formula) with a companion column F
# "Deprecated since?" that only held a value for deprecated rows.
#
# The new layout replaces the boolean with a textual "State" column
# (values "deprecated"/"active") and renames F to "Deprecation version",
# reusing column D's cell formatting (text number format) for the new F
# values/header since they are now plain data copied in the same shape as
# column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
$ws.Range("E1").Value = "State"
$ws.Range("F1").Value = "Deprecation version"

# --- Column E: boolean flag -> textual state -------------------------
$ws.Range("E2").Value = "deprecated"
$ws.Range("E3").Value = "deprecated"
$ws.Range("E4").Value = "deprecated"
$ws.Range("E5").Value = "active"
$ws.Range("E6").Value = "active"

# --- Column F: copy column D's number format onto the (re-purposed)
# header and data cells, same as was done by hand in the workbook this
# mirrors (copy/paste-format from D) -----------------------------------
$ws.Range("D1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$ws.Range("D2:D4").Copy()
$ws.Range("F2:F4").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Column widths: E keeps its width (content no longer needs
# bestFit auto-sizing), F grows to fit "Deprecation version" ----------
$ws.Columns("E").AutoFit()
$ws.Columns("F").ColumnWidth = 20.42

# --- Selection, matching where the edit left the cursor ---------------
$ws.Range("F5").Select()
